$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so that numeric-looking
# strings (e.g. "10.50", "0.00001030") keep their exact text representation
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.588.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.479.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9712'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '279.34'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3663'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3078'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.05'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06675'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.531'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.216'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9707'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001030'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.478.19'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05938'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.61'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.504'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.53'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.265'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.640.85'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.04'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.135'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.68%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.638.14'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.919'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8264'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.025'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07989'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.536'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.207'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05791'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.741'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9705'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02048'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.50'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.637'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1878'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5305'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.535'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.23'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.45'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5199'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.809'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06493'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9952'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.18%  '
